$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I50").Value = "ba"
$ws.Range("J50").Value = "Appreciation"

$ws.Range("I52").Value = "aa"
$ws.Range("J52").Value = "Agree/Accept"

$ws.Range("I56").Value = "sd"
$ws.Range("J56").Value = "Statement-non-opinion"

$ws.Range("I62").Value = "%"
$ws.Range("J62").Value = "Uninterpretable"

$ws.Range("I67").Value = "sd"
$ws.Range("J67").Value = "Statement-non-opinion"

$ws.Range("I72").Value = "aa"
$ws.Range("J72").Value = "Agree/Accept"

$ws.Range("I75").Value = "sd"
$ws.Range("J75").Value = "Statement-non-opinion"

$ws.Range("I77").Value = "aa"
$ws.Range("J77").Value = "Agree/Accept"

$ws.Range("I79").Value = "ba"
$ws.Range("J79").Value = "Appreciation"
